$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) are text-typed ("inlineStr") in the source
# file even though many values look numeric (e.g. "1.0000", "238.39").
# Force Text format on the whole touched range first so Excel does not
# silently coerce these into numbers (which would drop trailing zeros /
# collapse values like "1.0000" to "1").
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.561.99'
$ws.Range("E2").Value = '  +0.56%  '

$ws.Range("D3").Value = '1.881.95'
$ws.Range("E3").Value = '  -0.39%  '

$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  -0.35%  '

$ws.Range("D5").Value = '238.39'
$ws.Range("E5").Value = '  -0.01%  '

$ws.Range("D6").Value = '0.9992'
$ws.Range("E6").Value = '  -0.37%  '

$ws.Range("D7").Value = '0.4802'
$ws.Range("E7").Value = '  -0.59%  '

$ws.Range("D8").Value = '0.2820'
$ws.Range("E8").Value = '  -1.92%  '

$ws.Range("D9").Value = '0.06513'
$ws.Range("E9").Value = '  -1.39%  '

$ws.Range("D10").Value = '1.957.02'
$ws.Range("E10").Value = '  +3.37%  '

$ws.Range("D11").Value = '0.07474'
$ws.Range("E11").Value = '  +0.70%  '

$ws.Range("D12").Value = '16.56'
$ws.Range("E12").Value = '  -1.15%  '

$ws.Range("D13").Value = '5.087'
$ws.Range("E13").Value = '  -2.14%  '

$ws.Range("D14").Value = '87.86'
$ws.Range("E14").Value = '  -0.65%  '

$ws.Range("D15").Value = '0.6625'
$ws.Range("E15").Value = '  +0.47%  '

$ws.Range("D16").Value = '30.508.22'
$ws.Range("E16").Value = '  +0.40%  '

$ws.Range("B17").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C17").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D17").Value = '2.240.63'
$ws.Range("E17").Value = '  +4.98%  '

$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").Value = '13.25'
$ws.Range("E18").Value = '  -1.93%  '

$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '1.0000'
$ws.Range("E19").Value = '  -0.17%  '

$ws.Range("D20").Value = '0.000007555'
$ws.Range("E20").Value = '  -2.64%  '

$ws.Range("D21").Value = '228.40'
$ws.Range("E21").Value = '  +3.70%  '

$ws.Range("D22").Value = '0.9999'
$ws.Range("E22").Value = '  -0.41%  '

$ws.Range("D23").Value = '5.279'
$ws.Range("E23").Value = '  -1.77%  '

$ws.Range("D24").Value = '6.177'
$ws.Range("E24").Value = '  -0.35%  '

$ws.Range("D25").Value = '9.302'
$ws.Range("E25").Value = '  -1.22%  '

$ws.Range("D26").Value = '167.44'
$ws.Range("E26").Value = '  +1.56%  '

$ws.Range("D27").Value = '18.44'
$ws.Range("E27").Value = '  -2.73%  '

$ws.Range("D28").Value = '1.945'
$ws.Range("E28").Value = '  +0.06%  '

$ws.Range("D29").Value = '1.402'
$ws.Range("E29").Value = '  -3.94%  '

$ws.Range("D30").Value = '0.09742'
$ws.Range("E30").Value = '  +5.20%  '

$ws.Range("D31").Value = '4.332'
$ws.Range("E31").Value = '  +0.51%  '

$ws.Range("D32").Value = '4.005'
$ws.Range("E32").Value = '  -0.43%  '

$ws.Range("D33").Value = '0.05039'
$ws.Range("E33").Value = '  -0.84%  '

$ws.Range("D34").Value = '1.226'
$ws.Range("E34").Value = '  +6.49%  '

$ws.Range("D35").Value = '0.7511'
$ws.Range("E35").Value = '  -0.76%  '

$ws.Range("D36").Value = '2.708'
$ws.Range("E36").Value = '  +0.35%  '

$ws.Range("D37").Value = '0.01863'
$ws.Range("E37").Value = '  +0.12%  '

$ws.Range("D38").Value = '2.647'
$ws.Range("E38").Value = '  +0.19%  '

$ws.Range("D39").Value = '0.9109'
$ws.Range("E39").Value = '  -0.44%  '

$ws.Range("D40").Value = '2.070'
$ws.Range("E40").Value = '  -0.61%  '

$ws.Range("D41").Value = '106.50'
$ws.Range("E41").Value = '  -1.27%  '

$ws.Range("D42").Value = '0.4274'
$ws.Range("E42").Value = '  -1.33%  '

$ws.Range("D43").Value = '5.777'
$ws.Range("E43").Value = '  -3.05%  '

$ws.Range("D44").Value = '0.9991'
$ws.Range("E44").Value = '  -0.46%  '

$ws.Range("D45").Value = '7.358'
$ws.Range("E45").Value = '  -3.34%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '64.31'
$ws.Range("E46").Value = '  -1.25%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.1274'
$ws.Range("E47").Value = '  -3.93%  '

$ws.Range("D48").Value = '1.471'
$ws.Range("E48").Value = '  -7.45%  '

$ws.Range("D49").Value = '8.895'
$ws.Range("E49").Value = '  -1.13%  '

$ws.Range("D50").Value = '33.58'
$ws.Range("E50").Value = '  -3.19%  '

$ws.Range("D51").Value = '0.05650'
$ws.Range("E51").Value = '  -1.40%  '
